$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header
$ws.Range("B1").Value2 = "RangeDates"

# Add new column C
$ws.Range("C1").Value2 = "AloneDate"
$ws.Range("C2").Value2 = "2023-6-17,2023-6-18"

# Update existing value (long range string) last so it lands at the end
# of the shared strings table, matching the authored order
$ws.Range("B2").Value2 = "(2023-05-12, 2023-6-16), (2023-07-1, 2023-7-3)"

# Update selection to reflect the new active cell
$ws.Range("C2").Select()
